# Update the "Qmin" column (H) on the "Generators" sheet so that it equals
# the negative of the "Qmax" column (G) for each generator row (rows 2-54).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generators")

$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $qmax = $ws.Cells.Item($r, 7).Value2   # column G = Qmax
    $ws.Cells.Item($r, 8).Value2 = -$qmax  # column H = Qmin
}
